$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# Insert a new row above the current row 5 (the blank separator row) so the
# new MaxContinuousRetryNumber setting lands right after MaxInitRetryNumber.
$ws.Rows.Item(5).Insert()

# Reword the existing retry-related descriptions (content only; the shared
# string table index shift for these is handled automatically).
$ws.Range("C3").Value = "If > 0, the robot will retry the same transaction which failed with application exception. This is a local data retry. Orchestrator Queue Item retry are managed at the queue level. Must be integer"
$ws.Range("C4").Value = "If > 0 will retry the Initialisation state with a failed exception. Must be an integer."

# New setting row: MaxContinuousRetryNumber
$ws.Range("A5").Value = "MaxContinuousRetryNumber"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = "If > 0 will keep a record of consecutive failed exceptions of the Process state. When this number is reached, the application will fail. Must be an integer."

# Make Constants the active/selected sheet (was Workblocks before).
$ws.Activate()
$ws.Range("B6").Select()

Write-Output "done"
